$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated JobMatchIndex values (column C) resulting from the dataframe join fix
$matchIndexUpdates = @{
    2 = 1
    3 = 1
    4 = 2
    5 = 1
    8 = 3
    9 = 0
    11 = 2
    13 = 2
    14 = 0
    16 = 2
    17 = -4
    18 = 5
    19 = 0
    20 = -1
    22 = 0
    24 = 0
    25 = 1
    26 = 0
    27 = 2
    28 = -1
    29 = 1
    30 = 3
    31 = 3
    32 = 0
    34 = 2
    36 = -2
    37 = 0
    39 = -2
    40 = 0
    41 = -3
    42 = -1
    44 = 4
    45 = 1
    46 = 0
    47 = -5
    48 = -3
    49 = -1
    50 = -3
    51 = -1
    52 = 2
    55 = 0
    56 = -1
    57 = 0
    58 = 2
    59 = 1
    60 = -2
    61 = -1
    62 = -1
    63 = 0
    66 = 1
    69 = 5
    70 = 5
    71 = 0
    72 = 1
    74 = 0
    75 = 6
    76 = -1
    77 = 4
    78 = 0
    79 = -1
    80 = 0
    81 = 6
    82 = 0
    83 = 2
    84 = 0
    86 = 1
    87 = 2
    88 = 0
    89 = 2
    90 = -2
    91 = -2
    92 = -1
    93 = 0
    94 = 2
    95 = 1
    96 = 0
    98 = 2
    99 = -1
    100 = -2
    102 = 0
    103 = 1
    104 = 1
    105 = 0
    106 = 0
    107 = 3
    108 = 6
    109 = 5
    110 = -2
    111 = -2
    113 = -1
    114 = 4
    115 = -2
    116 = 2
    118 = 1
    119 = -1
    120 = 1
    121 = 1
    122 = 0
    123 = -1
    124 = 2
    125 = 0
    126 = 9
    127 = -1
    128 = -2
    129 = -3
    130 = 1
    131 = -1
    132 = 0
    133 = 0
    134 = 0
    135 = 0
    136 = 0
    137 = 0
    138 = -1
    139 = -1
    140 = -1
    141 = -1
    142 = 1
    145 = 2
    146 = -4
    148 = 0
    149 = 1
    151 = 0
    152 = -1
    153 = -1
    154 = 2
    155 = 0
    156 = 0
    158 = 0
    161 = 2
    162 = 0
    163 = 1
    165 = -1
    166 = -5
    168 = 2
    170 = 0
    171 = 1
    172 = 1
    173 = 0
    174 = 1
    176 = 2
    177 = 0
    178 = -1
    179 = 1
    181 = 4
    182 = 3
    183 = -4
    184 = -2
    185 = 4
    186 = 0
    188 = -1
    191 = 0
    192 = -1
    193 = -1
    194 = 2
    195 = -1
    197 = 2
    198 = 0
    199 = 1
    201 = 2
    202 = 0
    203 = 5
    204 = -1
    205 = -3
    206 = -1
    207 = -1
    208 = 1
    209 = 4
    210 = -1
}

foreach ($row in $matchIndexUpdates.Keys) {
    $ws.Cells.Item($row, 3).Value = $matchIndexUpdates[$row]
}
